$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (45488 -> 45489, i.e. 2024-07-15 -> 2024-07-16) for every data row (2..28).
for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45488) {
        $cell.Value = 45489
    }
}
